# Allow for more complex rules for data selection and add test function.
#
# The "Condition" column (B) previously contained a single repeated value "a"
# for every data row. This update differentiates the condition labels so
# more complex data-selection rules can be tested:
#   - most rows become "b"
#   - rows 3 and 6 become "ba"
#   - row 10 becomes "ab"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "b"
$ws.Range("B3").Value  = "ba"
$ws.Range("B4").Value  = "b"
$ws.Range("B5").Value  = "b"
$ws.Range("B6").Value  = "ba"
$ws.Range("B7").Value  = "b"
$ws.Range("B8").Value  = "b"
$ws.Range("B9").Value  = "b"
$ws.Range("B10").Value = "ab"
$ws.Range("B11").Value = "b"

# Reflect the author's final cell selection in the saved sheet view.
$ws.Range("A14").Select()
